$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 15 rework: the "estimated Amsteel II+ stretch" / 0.005 row is cleared out ---
# A15 previously held the shared string "estimated Amsteel II+ stretch" (style carried
# an applied border). Clear its content and restyle it to match the plain "H column"
# look used throughout this table (same formatting as H5/H6/H8/etc.).
$ws.Range("A15").ClearContents()
$ws.Range("H5").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# B15 (the 0.5% value) and H15 are fully removed - not just blanked.
$ws.Range("B15").Clear()
$ws.Range("H15").Clear()

# C15:G15 become blank cells carrying the workbook's plain default look (same as
# columns A/B/H's own default style). Build that style on a scratch cell first
# (column J already carries it), copy its format across, then restore the scratch
# cell to its untouched state.
$ws.Range("J1").Font.Size = 11
$ws.Range("J1").Copy()
$ws.Range("C15:G15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("J1").Clear()

# --- Selection moved from C21 to A3:H3 ---
$ws.Range("A3:H3").Select() | Out-Null
